# Fix typo in Cheat sheet with regex (#25215)
#
# Slide 1 contains a "regex cheatsheet" table (shape "Table 84") whose
# last data row documents the regex that matches every column except
# "Species". The regex example cell had an accidental doubled leading
# single-quote: "''^(?!Species$).*'" instead of "'^(?!Species$).*'".
# This fixes the typo, leaving all text formatting untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$fixed = $false
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Table 84" -and $shp.HasTable) {
        $tbl = $shp.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cellRange = $tbl.Cell($r, $c).Shape.TextFrame.TextRange
                if ($cellRange.Text -eq "''^(?!Species`$).*'") {
                    $cellRange.Text = "'^(?!Species`$).*'"
                    $fixed = $true
                }
            }
        }
    }
}

if (-not $fixed) {
    throw "Could not locate the regex typo cell to fix"
}
